$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price values are stored as text (matches original inlineStr type)
$textCells = @('D9', 'D11', 'D20', 'D26', 'D27', 'D30', 'D31', 'D41', 'D42', 'D43', 'D47')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values from the refreshed cryptos feed
$ws.Range('D2').Value = '25.302.46'
$ws.Range('E2').Value = '  -2.49%  '
$ws.Range('D3').Value = '1.566.65'
$ws.Range('E4').Value = '  -0.37%  '
$ws.Range('E5').Value = '  -2.88%  '
$ws.Range('E6').Value = '  -0.35%  '
$ws.Range('E7').Value = '  -5.01%  '
$ws.Range('E8').Value = '  -1.52%  '
$ws.Range('D9').Value = '0.243'
$ws.Range('E9').Value = '  -2.88%  '
$ws.Range('E10').Value = '  -2.44%  '
$ws.Range('D11').Value = '0.0781'
$ws.Range('E11').Value = '  -0.85%  '
$ws.Range('D12').Value = '1.784.92'
$ws.Range('E12').Value = '  -3.54%  '
$ws.Range('D13').Value = '1.565.09'
$ws.Range('E13').Value = '  -3.87%  '
$ws.Range('E14').Value = '  -3.93%  '
$ws.Range('E15').Value = '  -3.43%  '
$ws.Range('D16').Value = '25.301.67'
$ws.Range('E16').Value = '  -2.42%  '
$ws.Range('E17').Value = '  -2.90%  '
$ws.Range('D18').Value = '0.0₃0712'
$ws.Range('E18').Value = '  -3.30%  '
$ws.Range('E19').Value = '  -0.29%  '
$ws.Range('D20').Value = '185.54'
$ws.Range('E20').Value = '  -3.19%  '
$ws.Range('E21').Value = '  -2.44%  '
$ws.Range('E22').Value = '  -2.98%  '
$ws.Range('E23').Value = '  -3.16%  '
$ws.Range('E24').Value = '  -0.35%  '
$ws.Range('E25').Value = '  -3.94%  '
$ws.Range('D26').Value = '140.35'
$ws.Range('E26').Value = '  -2.26%  '
$ws.Range('D27').Value = '1.64'
$ws.Range('E27').Value = '  -6.92%  '
$ws.Range('E28').Value = '  -3.81%  '
$ws.Range('E29').Value = '  -2.17%  '
$ws.Range('D30').Value = '1.15'
$ws.Range('E30').Value = '  -6.14%  '
$ws.Range('D31').Value = '0.0464'
$ws.Range('E31').Value = '  -3.70%  '
$ws.Range('E32').Value = '  -2.79%  '
$ws.Range('E33').Value = '  -3.39%  '
$ws.Range('E34').Value = '  -2.00%  '
$ws.Range('E35').Value = '  -3.43%  '
$ws.Range('D36').Value = '1.090.37'
$ws.Range('E36').Value = '  -2.76%  '
$ws.Range('E37').Value = '  -0.72%  '
$ws.Range('E38').Value = '  -4.84%  '
$ws.Range('E39').Value = '  -2.62%  '
$ws.Range('E40').Value = '  -4.44%  '
$ws.Range('D41').Value = '0.769'
$ws.Range('E41').Value = '  -9.09%  '
$ws.Range('D42').Value = '0.784'
$ws.Range('E42').Value = '  +2.09%  '
$ws.Range('D43').Value = '93.43'
$ws.Range('E43').Value = '  -4.47%  '
$ws.Range('E44').Value = '  -2.02%  '
$ws.Range('D45').Value = '1.699.14'
$ws.Range('E45').Value = '  -3.51%  '
$ws.Range('E46').Value = '  -2.65%  '
$ws.Range('D47').Value = '52.68'
$ws.Range('E47').Value = '  -3.36%  '
$ws.Range('E48').Value = '  -4.63%  '
$ws.Range('E49').Value = '  -2.30%  '
$ws.Range('E50').Value = '  -1.67%  '
$ws.Range('E51').Value = '  -0.50%  '
